$d = $word.ActiveDocument

# Update the report date shown next to the author name
$d.Content.Find.Execute("26 Aug 2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "12 Sep 2024", 2)

# Update the git revision used to generate the report
$d.Content.Find.Execute("1ba5e7f", $true, $false, $false, $false, $false,
                         $true, 1, $false, "dc28f0e", 2)
